$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (losing the original text formatting,
# e.g. trailing zeros), matching the source file where these are text cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.505.28"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "2.906.42"
$ws.Range("E3").Value = "  -3.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "527.62"
$ws.Range("E5").Value = "  -5.39%  "
$ws.Range("D6").Value = "141.97"
$ws.Range("E6").Value = "  -7.97%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").Value = "2.907.27"
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -5.32%  "
$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  -8.81%  "
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").Value = "3.413.43"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "60.649.60"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "22.60"
$ws.Range("E16").Value = "  -6.08%  "
$ws.Range("D17").Value = "2.904.97"
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("E18").Value = "  -6.66%  "
$ws.Range("D19").Value = "4.89"
$ws.Range("E19").Value = "  -4.31%  "
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -3.27%  "
$ws.Range("D21").Value = "359.21"
$ws.Range("E21").Value = "  -9.33%  "
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "5.67"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "63.35"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").Value = "3.029.25"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").Value = "0.448"
$ws.Range("E27").Value = "  -4.05%  "
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  -5.42%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "0.0₃0857"
$ws.Range("E30").Value = "  -12.29%  "
$ws.Range("D31").Value = "7.61"
$ws.Range("E31").Value = "  -12.55%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "1.66"
$ws.Range("E33").Value = "  -5.13%  "
$ws.Range("D34").Value = "19.64"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").Value = "153.52"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").Value = "4.30"
$ws.Range("E36").Value = "  -9.02%  "
$ws.Range("D37").Value = "5.52"
$ws.Range("E37").Value = "  -8.70%  "
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -9.50%  "
$ws.Range("E39").Value = "  -8.36%  "
$ws.Range("D40").Value = "37.87"
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").Value = "2.321.25"
$ws.Range("E41").Value = "  -8.23%  "
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  -8.13%  "
$ws.Range("D43").Value = "3.66"
$ws.Range("E43").Value = "  -6.70%  "
$ws.Range("D44").Value = "0.641"
$ws.Range("E44").Value = "  -4.13%  "
$ws.Range("D45").Value = "20.71"
$ws.Range("E45").Value = "  -9.09%  "
$ws.Range("D46").Value = "0.0567"
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "4.83"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "10.33"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0232"
$ws.Range("E50").Value = "  -6.61%  "
$ws.Range("D51").Value = "0.0920"
$ws.Range("E51").Value = "  -2.66%  "
